$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row contents ("Key" / "Value") without shifting the rows below
$ws.Range("A1:B1").ClearContents()

# Update selection to match the post-edit state
$ws.Range("C6").Select()
